# Apply the "Penalty Reward System" forecast refresh to the workbook.
# Sheet "Forecast Comparison": shift Week_Start_Date (col B) forward by one
# week and overwrite MyForecast (col D) with the new predicted values for
# rows 2-17. Sheet "Summary": refresh the derived statistics accordingly.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date values (col B), rows 2-17
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast values (col D), rows 2-17
$newForecast = @(53, 64, 81, 87, 76, 66, 73, 93, 102, 92, 78, 83, 71, 68, 68, 89)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $dateCell = $wsForecast.Cells.Item($row, 2)
    # Keep the date stored as plain text (matches the original inline-string layout)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDates[$i]
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

# Refresh the Summary sheet stats (also plain text cells)
$wsSummary.Range("B2").NumberFormat  = "@"
$wsSummary.Range("B2").Value  = "2024-07-14 to 2025-01-05"
$wsSummary.Range("B4").Value  = "97"
$wsSummary.Range("B5").Value  = "23"
$wsSummary.Range("B6").Value  = "10"
$wsSummary.Range("B7").Value  = "27"
$wsSummary.Range("B8").Value  = "595 units"
$wsSummary.Range("B9").Value  = "1244"
$wsSummary.Range("B10").Value = "593"
$wsSummary.Range("B11").Value = "285"
$wsSummary.Range("B12").Value = "102"
